$d = $word.ActiveDocument

# --- 1) Fix the grammar in the "Game: Have you?" closing sentence ---
# Original: " If the other person haven’t have this get 1 points."
# Target:   " If the other person hasn’t had this get 1 point."
$d.Content.Find.Execute("haven’t have this get 1 points.", $true, $false, $false, $false, $false, $true, 1, $false, "hasn’t had this get 1 point.", 2)

# --- 2) Add a new "cours of the day" section at the end of the document ---
# Insert two new paragraphs right before the final (already existing) empty paragraph:
#   - a blank "Text" styled paragraph
#   - a "Categorie" styled paragraph with the text "CV / Resume"
$last = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $last.Range
$lastRange.InsertParagraphBefore()
$lastRange.InsertParagraphBefore()

$newCategorie = $d.Paragraphs($d.Paragraphs.Count - 1)
$newCategorie.Style = "Categorie"
$newCategorie.Range.Text = "CV / Resume"
